# Generate Report for Handoff
#
# Refresh the handoff report timestamps and mark the 6 "awaiting handoff"
# rows (per locale sheet) with the "ht" (handoff) priority now that a new
# handoff xliff set was generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows whose handoff batch was just regenerated.
$rows = @(7, 8, 9, 10, 11, 13)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" (column G)
    $wsOverview.Range("G$r").Value = "2016-08-22 20:19:55"

    # zh-cn sheet: "Priority" (column E) + "Latest Handoff Datetime" (column H)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-22 20:19:50"

    # de-de sheet: "Priority" (column E) + "Latest Handoff Datetime" (column H)
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-22 20:19:55"
}
